# Production doc builds - 2021/04/08 00:59:32 UTC
#
# The only substantive, author-driven content change in this revision is on
# slide 1: the second "Availability Zone 1" label (shape "Rectangle 50",
# the dashed rectangle around the right-hand VPC availability zone) is
# updated to read "Availability Zone 2". Editing it in place via
# TextRange.Characters (rather than replacing the whole TextRange.Text)
# mirrors how PowerPoint itself splits the run when only the trailing
# character is retyped, so the untouched "Availability Zone " prefix keeps
# its original run formatting and only the changed "1" -> "2" character
# becomes (part of) a new run.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("Rectangle 50")

$tr = $shape.TextFrame.TextRange
$len = $tr.Length
$lastChar = $tr.Characters($len, 1)
$lastChar.Text = "2"
